$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column I (9) so the old I/J/K (dv/dt, dv/dt m/s2, pred move)
# shift right to J/K/L, making room for the new "v, m/s" column at I.
$ws.Columns.Item(9).Insert()

# --- Header row -----------------------------------------------------------
$ws.Range("B1").Value = "t, sec"
$ws.Range("C1").Value = "v, mph"
$ws.Range("D1").Value = "px, un"
$ws.Range("E1").Value = "py, un"
$ws.Range("F1").Value = "move, un"
$ws.Range("G1").Value = "dt, sec"
$ws.Range("H1").Value = "move/dt, un/s"
$ws.Range("I1").Value = "v, m/s"
$ws.Range("J1").Value = "dv, m/s"
$ws.Range("K1").Value = "dv/dt, m/s2"
$ws.Range("L1").Value = "pred move, m"

# --- Column H: move/dt, un/s = F/G (was dv = C-Cprev) ----------------------
$ws.Range("H3:H6").Formula = "=F3/G3"

# --- Column I: v, m/s = C*1609/3600 (new column, includes row 2) -----------
$ws.Range("I2:I6").Formula = "=C2*1609/3600"

# --- Column J: dv, m/s = I - Iprev ------------------------------------------
$ws.Range("J3:J6").Formula = "=I3-I2"

# --- Column K: dv/dt, m/s2 = J/G --------------------------------------------
$ws.Range("K3:K6").Formula = "=J3/G3"

# --- Column L: pred move, m = Iprev*G + K/2*G^2 -----------------------------
$ws.Range("L3:L6").Formula = "=I2*G3+K3/2*G3^2"

# --- column widths / selection ----------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 13.85546875
$ws.Columns.Item(9).ColumnWidth = 13.85546875
$ws.Columns.Item(11).ColumnWidth = 14.85546875
$ws.Columns.Item(12).ColumnWidth = 12

$ws.Range("O10").Select()
